# "Marked all tests for execution Including upload document test"
#
# The "Test Cases" sheet tracks Sanity Runmode (col D: Y/N) and the last
# Results (col E: PASS/FAIL) for each test case. Marking a test case for
# (re-)execution means flipping its Sanity Runmode to "Y" and clearing any
# stale result so it shows as not-yet-run.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Mark the two not-yet-marked test cases (rows 2 and 3) for execution.
$ws.Range("D2").Value = "Y"
$ws.Range("D3").Value = "Y"

# Clear stale PASS/FAIL results now that these rows are queued to run again.
$ws.Range("E3").ClearContents()
$ws.Range("E4").ClearContents()

# Leave the selection where the author left it.
$ws.Range("D3").Select()
